# Bug fix for add_content_slide: each "content" slide (slides 2-5) had an
# extra decorative "Rectangle" shape (the gradient header bar) that should
# not have been added by the generator. Remove it (and any duplicate of it)
# from every content slide, then renumber the remaining shapes' names so
# they stay sequential (TextBox N / Picture N), matching what the slide
# generator would have produced without the stray rectangle(s).

$p = $ppt.ActivePresentation

for ($slideIdx = 2; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)

    # Remove every stray "Rectangle*" shape (the gradient bar bug).
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -like "Rectangle*") {
            $sh.Delete()
        }
    }

    # Renumber the surviving shapes so names stay sequential again.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -like "TextBox*") {
            $sh.Name = "TextBox " + $i
        } elseif ($sh.Name -like "Picture*") {
            $sh.Name = "Picture " + $i
        }
    }
}
